$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking price/percentage cells so Excel
# does not coerce them into numbers (losing literal formatting such as
# trailing zeros, e.g. "9.160" or "0.10%").
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '330.58'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.46%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '41.46'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '0.27%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.698'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.15%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08427'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '4.40%'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.14%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '4.495'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.43%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.985'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-3.41%'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '1.06%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9269'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '0.64%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1261'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.58%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1982'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '1.61%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09511'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '1.44%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03968'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '8.02%'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.97%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001303'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.45%'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04425'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.06%'
$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006114'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.36%'
$ws.Range("B19").Value = 'LEO'
$ws.Range("C19").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.435'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.59%'
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3511'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.84%'
$ws.Range("B21").Value = 'MCDex'
$ws.Range("C21").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.160'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '10.59%'
$ws.Range("B22").Value = 'ProBitToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.1364'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-3.71%'
$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.2512'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-5.25%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001246'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.93%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004398'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.55%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-3.97%'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '0.07%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02832'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '0.04%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05520'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.83%'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '4.17%'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '1.45%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.008968'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-9.81%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002083'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-1.32%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01097'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-7.71%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00007311'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '8.38%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.10%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003221'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '7.72%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002281'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.06%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002103'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.10%'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.10%'
